$wb = $excel.ActiveWorkbook

# ALC (sheet1) row 31
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2725.3
$ws.Range("I31").Value = 1124.8334
$ws.Range("J31").Value = 5126
$ws.Range("K31").Value = 3374.5002
$ws.Range("L31").Value = 15378
$ws.Range("M31").Value = -3144.5002
$ws.Range("N31").Value = -15838

# ALC (sheet1) row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 107.375
$ws.Range("I39").Value = 78.166664
$ws.Range("J39").Value = 195
$ws.Range("K39").Value = 234.499992
$ws.Range("L39").Value = 585
$ws.Range("M39").Value = 61.50000800000001
$ws.Range("N39").Value = -1177

# ALC (sheet1) row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2060.4167
$ws.Range("I40").Value = 1638.8889
$ws.Range("J40").Value = 2313.3333
$ws.Range("K40").Value = 1638.8889
$ws.Range("L40").Value = 2313.3333
$ws.Range("M40").Value = -1463.8889
$ws.Range("N40").Value = -2663.3333

# ALC (sheet1) row 42
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 64.28570999999999
$ws.Range("I42").Value = 50.2
$ws.Range("J42").Value = 99.5
$ws.Range("K42").Value = 150.6
$ws.Range("L42").Value = 298.5
$ws.Range("M42").Value = 79.39999999999998
$ws.Range("N42").Value = -758.5

# ALC (sheet1) row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6123.1665
$ws.Range("I116").Value = 10549.75
$ws.Range("J116").Value = 3909.875
$ws.Range("K116").Value = 10549.75
$ws.Range("L116").Value = 3909.875
$ws.Range("M116").Value = -7107.75
$ws.Range("N116").Value = -10793.875

# ALC (sheet1) row 118
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 2712.923
$ws.Range("I118").Value = 1964.2858
$ws.Range("J118").Value = 3586.3333
$ws.Range("K118").Value = 5892.857400000001
$ws.Range("L118").Value = 10758.9999
$ws.Range("M118").Value = -4235.857400000001
$ws.Range("N118").Value = -14072.9999

# ALC (sheet1) row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2569.443
$ws.Range("I129").Value = 497.8889
$ws.Range("J129").Value = 2805.443
$ws.Range("K129").Value = 1493.6667
$ws.Range("L129").Value = 8416.329000000002
$ws.Range("M129").Value = 3506.3333
$ws.Range("N129").Value = -18416.329

# ALC (sheet1) row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5558570.5
$ws.Range("I138").Value = 1585.5454
$ws.Range("J138").Value = 10873947
$ws.Range("K138").Value = 4756.6362
$ws.Range("L138").Value = 32621841
$ws.Range("M138").Value = 383.3638000000001
$ws.Range("N138").Value = -32632121

# ARM (sheet2) row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6277.421
$ws.Range("I32").Value = 6780.6
$ws.Range("K32").Value = 6780.6
$ws.Range("M32").Value = -6493.6

# ARM (sheet2) row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2663.2964
$ws.Range("I45").Value = 3582.4
$ws.Range("J45").Value = 2122.647
$ws.Range("K45").Value = 3582.4
$ws.Range("L45").Value = 2122.647
$ws.Range("M45").Value = -3205.4
$ws.Range("N45").Value = -2876.647

# ARM (sheet2) row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1211.5
$ws.Range("I110").Value = 1229.3077
$ws.Range("J110").Value = 980
$ws.Range("K110").Value = 1229.3077
$ws.Range("L110").Value = 980
$ws.Range("M110").Value = 815.6922999999999
$ws.Range("N110").Value = -5070

# BSM (sheet3) row 75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 32503.385
$ws.Range("J75").Value = 50220.57
$ws.Range("L75").Value = 50220.57
$ws.Range("N75").Value = -52092.57

# BSM (sheet3) row 78
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H78").Value = 32503.385
$ws.Range("J78").Value = 50220.57
$ws.Range("L78").Value = 150661.71
$ws.Range("N78").Value = -160021.71

# BSM (sheet3) row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 916.6923
$ws.Range("I94").Value = 854.2727
$ws.Range("J94").Value = 1260
$ws.Range("K94").Value = 854.2727
$ws.Range("L94").Value = 1260
$ws.Range("M94").Value = -403.2727
$ws.Range("N94").Value = -2162

# CUL (sheet5) row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1157609.4
$ws.Range("I2").Value = 497.5
$ws.Range("J2").Value = 1736165.2
$ws.Range("K2").Value = 2985
$ws.Range("L2").Value = 10416991.2
$ws.Range("M2").Value = -2872
$ws.Range("N2").Value = -10417217.2

# CUL (sheet5) row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 552.0417
$ws.Range("I14").Value = 552.0417
$ws.Range("K14").Value = 1656.1251
$ws.Range("M14").Value = -1483.1251

# CUL (sheet5) row 35
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 3000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 3000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 9000
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -9576

# CUL (sheet5) row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 115.388885
$ws.Range("I38").Value = 15.571428
$ws.Range("J38").Value = 178.90909
$ws.Range("K38").Value = 46.714284
$ws.Range("L38").Value = 536.72727
$ws.Range("M38").Value = 300.285716
$ws.Range("N38").Value = -1230.72727

# CUL (sheet5) row 42
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 7775
$ws.Range("J42").Value = 7775
$ws.Range("L42").Value = 23325
$ws.Range("N42").Value = -24393

# CUL (sheet5) row 64
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2714.2856
$ws.Range("J64").Value = 3400
$ws.Range("L64").Value = 10200
$ws.Range("N64").Value = -10740

# CUL (sheet5) row 67
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 2714.2856
$ws.Range("J67").Value = 3400
$ws.Range("L67").Value = 10200
$ws.Range("N67").Value = -12072

# CUL (sheet5) row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 867.6076
$ws.Range("J68").Value = 955.2245
$ws.Range("L68").Value = 2865.6735
$ws.Range("N68").Value = -4487.6735

# CUL (sheet5) row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 867.6076
$ws.Range("J71").Value = 955.2245
$ws.Range("L71").Value = 8597.020500000001
$ws.Range("N71").Value = -16709.0205

# CUL (sheet5) row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2900
$ws.Range("I80").Value = 1400
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 4200
$ws.Range("L80").Value = 10500
$ws.Range("M80").Value = -3264
$ws.Range("N80").Value = -12372

# CUL (sheet5) row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 2900
$ws.Range("I83").Value = 1400
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 12600
$ws.Range("L83").Value = 31500
$ws.Range("M83").Value = -7920
$ws.Range("N83").Value = -40860

# CUL (sheet5) row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 620.5441
$ws.Range("I107").Value = 286.67743
$ws.Range("J107").Value = 900.27026
$ws.Range("K107").Value = 860.0322900000001
$ws.Range("L107").Value = 2700.81078
$ws.Range("M107").Value = 1059.96771
$ws.Range("N107").Value = -6540.81078

# GSM (sheet6) row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 38.307693
$ws.Range("I2").Value = 48.57143
$ws.Range("J2").Value = 26.333334
$ws.Range("K2").Value = 48.57143
$ws.Range("L2").Value = 26.333334
$ws.Range("M2").Value = 64.42857000000001
$ws.Range("N2").Value = -252.333334

# GSM (sheet6) row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1251
$ws.Range("I113").Value = 1279.1875
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 1279.1875
$ws.Range("L113").Value = 800
$ws.Range("M113").Value = 890.8125
$ws.Range("N113").Value = -5140

# LTW (sheet7) row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2089.4443
$ws.Range("J61").Value = 1929.2858
$ws.Range("L61").Value = 1929.2858
$ws.Range("N61").Value = -2333.2858

# LTW (sheet7) row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1930.7646
$ws.Range("I82").Value = 1765.375
$ws.Range("J82").Value = 2077.7778
$ws.Range("K82").Value = 1765.375
$ws.Range("L82").Value = 2077.7778
$ws.Range("M82").Value = -1404.375
$ws.Range("N82").Value = -2799.7778

# LTW (sheet7) row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1930.7646
$ws.Range("I85").Value = 1765.375
$ws.Range("J85").Value = 2077.7778
$ws.Range("K85").Value = 1765.375
$ws.Range("L85").Value = 2077.7778
$ws.Range("M85").Value = -517.375
$ws.Range("N85").Value = -4573.7778

# LTW (sheet7) row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2089.4443
$ws.Range("J113").Value = 1929.2858
$ws.Range("L113").Value = 1929.2858
$ws.Range("N113").Value = -6269.2858

# WVR (sheet8) row 58
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 22933.334
$ws.Range("I58").Value = 21800
$ws.Range("J58").Value = 23500
$ws.Range("K58").Value = 21800
$ws.Range("L58").Value = 23500
$ws.Range("M58").Value = -21492
$ws.Range("N58").Value = -24116

# WVR (sheet8) row 75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 28381.111
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 28381.111
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 28381.111
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -30253.111

# WVR (sheet8) row 78
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 28381.111
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 28381.111
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 85143.333
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -94503.333
